$p = $ppt.ActivePresentation

# --- Slide 7: "2. Tight Coupled Code " (Database System Abstraction) -> "3. Tight Coupled Code "
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Title.TextFrame.TextRange
$tr7.InsertBefore("3") | Out-Null
$c7 = $tr7.Characters(1, 1)
$c7.Font.Size = 32
$c7.LanguageID = "bg-BG"
$tr7.Characters(2, 1).Text = ""

# --- Slide 8: "2. Loose Coupled Code " (Database System Abstraction) -> "3. Loose Coupled Code "
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Title.TextFrame.TextRange
$tr8.InsertBefore("3") | Out-Null
$c8 = $tr8.Characters(1, 1)
$c8.Font.Size = 32
$c8.LanguageID = "bg-BG"
$tr8.Characters(2, 1).Text = ""

# --- Slide 9: "2. Tight Coupled Code " (Client Notifier) -> "4. Tight Coupled Code "
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Title.TextFrame.TextRange
$tr9.InsertBefore("4") | Out-Null
$c9 = $tr9.Characters(1, 1)
$c9.Font.Size = 32
$c9.LanguageID = "bg-BG"
$tr9.Characters(2, 1).Text = ""

# --- Slide 10: "2. " + "Loose Coupled Code " (Client Notifier) -> "4. " + "Loose Coupled Code "
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Title.TextFrame.TextRange
$tr10.InsertBefore("4") | Out-Null
$c10 = $tr10.Characters(1, 1)
$c10.LanguageID = "bg-BG"
$tr10.Characters(2, 1).Text = ""
